$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 46038
$ws.Cells.Item(2, 2).Value = 12129.5764479979
$ws.Cells.Item(2, 3).Value = 10795.5808293969
$ws.Cells.Item(2, 4).Value = 15603.86
$ws.Cells.Item(2, 5).Value = 7917.26698617865
$ws.Cells.Item(2, 6).Value = 129.541158982316

$ws.Cells.Item(3, 1).Value = 46039
$ws.Cells.Item(3, 2).Value = 4664.54780201638
$ws.Cells.Item(3, 3).Value = 7300.10283619102
$ws.Cells.Item(3, 4).Value = 12075.86
$ws.Cells.Item(3, 5).Value = 7781.48427831776
$ws.Cells.Item(3, 6).Value = 125.238629771199

$ws.Cells.Item(4, 1).Value = 46040
$ws.Cells.Item(4, 2).Value = 4582.65855493572
$ws.Cells.Item(4, 3).Value = 7251.59088778583
$ws.Cells.Item(4, 4).Value = 12075.86
$ws.Cells.Item(4, 5).Value = 7780.10217081282
$ws.Cells.Item(4, 6).Value = 123.159710774944

$ws.Cells.Item(5, 1).Value = 46041
$ws.Cells.Item(5, 2).Value = 11846.49268181
$ws.Cells.Item(5, 3).Value = 11224.223643579
$ws.Cells.Item(5, 4).Value = 12075.86
$ws.Cells.Item(5, 5).Value = 8000.19718569214
$ws.Cells.Item(5, 6).Value = 297.85670121963

$ws.Cells.Item(6, 1).Value = 46042
$ws.Cells.Item(6, 2).Value = 12204.6399876786
$ws.Cells.Item(6, 3).Value = 11833.5823255322
$ws.Cells.Item(6, 4).Value = 12075.86
$ws.Cells.Item(6, 5).Value = 8240.15998098017
$ws.Cells.Item(6, 6).Value = 333.245096104681

$ws.Cells.Item(7, 1).Value = 46043
$ws.Cells.Item(7, 2).Value = 12733.0845017406
$ws.Cells.Item(7, 3).Value = 12151.2062012121
$ws.Cells.Item(7, 4).Value = 12075.86
$ws.Cells.Item(7, 5).Value = 8657.59313716841
$ws.Cells.Item(7, 6).Value = 363.87247243252

$ws.Cells.Item(8, 1).Value = 46044
$ws.Cells.Item(8, 2).Value = 12733.0845017406
$ws.Cells.Item(8, 3).Value = 12174.6151583264
$ws.Cells.Item(8, 4).Value = 12075.86
$ws.Cells.Item(8, 5).Value = 8657.59313716841
$ws.Cells.Item(8, 6).Value = 364.847845645616

$ws.Cells.Item(9, 1).Value = 46045
$ws.Cells.Item(9, 2).Value = 12733.0845017406
$ws.Cells.Item(9, 3).Value = 11594.2043814328
$ws.Cells.Item(9, 4).Value = 12075.86
$ws.Cells.Item(9, 5).Value = 8657.59313716841
$ws.Cells.Item(9, 6).Value = 340.66406327505

$ws.Cells.Item(10, 1).Value = 46046
$ws.Cells.Item(10, 2).Value = 5107.54318705847
$ws.Cells.Item(10, 3).Value = 8564.77604817398
$ws.Cells.Item(10, 4).Value = 12075.86
$ws.Cells.Item(10, 5).Value = 8264.44717782362
$ws.Cells.Item(10, 6).Value = 198.056801083233

$ws.Cells.Item(11, 1).Value = 46047
$ws.Cells.Item(11, 2).Value = 5000.01932310789
$ws.Cells.Item(11, 3).Value = 8740.20859929753
$ws.Cells.Item(11, 4).Value = 12075.86
$ws.Cells.Item(11, 5).Value = 8256.70262722195
$ws.Cells.Item(11, 6).Value = 205.043801104979

$ws.Cells.Item(12, 1).Value = 46048
$ws.Cells.Item(12, 2).Value = 12315.9682835607
$ws.Cells.Item(12, 3).Value = 12491.7285803074
$ws.Cells.Item(12, 4).Value = 12075.86
$ws.Cells.Item(12, 5).Value = 8314.28494857018
$ws.Cells.Item(12, 6).Value = 363.756397036568

$ws.Cells.Item(13, 1).Value = 46049
$ws.Cells.Item(13, 2).Value = 12315.9682835607
$ws.Cells.Item(13, 3).Value = 12369.7506383785
$ws.Cells.Item(13, 4).Value = 12075.86
$ws.Cells.Item(13, 5).Value = 8314.28494857018
$ws.Cells.Item(13, 6).Value = 358.673982789529

$ws.Cells.Item(14, 1).Value = 46050
$ws.Cells.Item(14, 2).Value = 12315.9682835607
$ws.Cells.Item(14, 3).Value = 11931.492762538
$ws.Cells.Item(14, 4).Value = 12075.86
$ws.Cells.Item(14, 5).Value = 8314.28494857018
$ws.Cells.Item(14, 6).Value = 340.413237962841

$ws.Cells.Item(15, 1).Value = 46051
$ws.Cells.Item(15, 2).Value = 12315.9682835607
$ws.Cells.Item(15, 3).Value = 12739.6794735367
$ws.Cells.Item(15, 4).Value = 12075.86
$ws.Cells.Item(15, 5).Value = 8314.28494857018
$ws.Cells.Item(15, 6).Value = 374.087684254452

Write-Host "done"